# Update column G ("K") values for rows 2-12 on Sheet1
# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 3
    3  = 1
    4  = 1
    5  = 3
    6  = 2
    7  = 2
    8  = 1
    9  = 2
    10 = 0
    11 = 2
    12 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
